# Apply updated crypto price/volume data per GitHub Actions scrape update
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Ensure price/volume cells remain stored as Text (not auto-converted to numbers),
# matching the original string cell semantics, before assigning their new values.
$textCells = @("D2", "E2", "D3", "E3", "D4", "E4", "D5", "D6", "E6", "D7", "E7", "D8", "E8", "D9", "E9", "D10", "E10", "D11", "D12", "E12", "D13", "E13", "D14", "E14", "D15", "E15", "D16", "E16", "D17", "E17", "D18", "E18", "D19", "E19", "D20", "E20", "D21", "E21", "D22", "E22", "D23", "E23", "D24", "E24", "D25", "E25", "D26", "E26", "D27", "E27", "D28", "E28", "D29", "E29", "D30", "E30", "D31", "E31", "D32", "E32", "D33", "E33", "D34", "E34", "D35", "E35", "D37", "E37", "D38", "E38", "D39", "E39", "D40", "E40", "D41", "E41", "D42", "E42", "D43", "E43", "D44", "E44", "D45", "D46", "E46", "D47", "D48", "E48", "D49", "E49", "D50", "E50", "D51", "E51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '23.175.66'
$ws.Range('E2').Value = '  +0.38%  '
$ws.Range('D3').Value = '1.599.98'
$ws.Range('E3').Value = '  -0.02%  '
$ws.Range('D4').Value = '1.002'
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').Value = '1.001'
$ws.Range('D6').Value = '302.92'
$ws.Range('E6').Value = '  +0.60%  '
$ws.Range('D7').Value = '0.3783'
$ws.Range('E7').Value = '  +0.10%  '
$ws.Range('D8').Value = '52.19'
$ws.Range('E8').Value = '  +4.71%  '
$ws.Range('D9').Value = '0.3611'
$ws.Range('E9').Value = '  -1.01%  '
$ws.Range('D10').Value = '1.263'
$ws.Range('E10').Value = '  -0.37%  '
$ws.Range('D11').Value = '1.002'
$ws.Range('D12').Value = '0.08120'
$ws.Range('E12').Value = '  -0.49%  '
$ws.Range('D13').Value = '22.61'
$ws.Range('E13').Value = '  -1.83%  '
$ws.Range('D14').Value = '6.571'
$ws.Range('E14').Value = '  -0.25%  '
$ws.Range('D15').Value = '7.383'
$ws.Range('E15').Value = '  +0.28%  '
$ws.Range('D16').Value = '0.00001246'
$ws.Range('E16').Value = '  -0.96%  '
$ws.Range('D17').Value = '1.597.85'
$ws.Range('E17').Value = '  -0.09%  '
$ws.Range('D18').Value = '93.89'
$ws.Range('E18').Value = '  +2.50%  '
$ws.Range('D19').Value = '0.06895'
$ws.Range('E19').Value = '  +0.52%  '
$ws.Range('D20').Value = '18.03'
$ws.Range('E20').Value = '  -1.58%  '
$ws.Range('D21').Value = '6.532'
$ws.Range('E21').Value = '  -0.43%  '
$ws.Range('D22').Value = '1.001'
$ws.Range('E22').Value = '  +0.00%  '
$ws.Range('D23').Value = '12.94'
$ws.Range('E23').Value = '  -0.18%  '
$ws.Range('D24').Value = '23.181.34'
$ws.Range('E24').Value = '  +0.41%  '
$ws.Range('D25').Value = '2.380'
$ws.Range('E25').Value = '  +1.65%  '
$ws.Range('D26').Value = '2.991'
$ws.Range('E26').Value = '  +10.20%  '
$ws.Range('D27').Value = '21.21'
$ws.Range('E27').Value = '  +0.41%  '
$ws.Range('D28').Value = '149.29'
$ws.Range('E28').Value = '  -0.70%  '
$ws.Range('D29').Value = '5.268'
$ws.Range('E29').Value = '  -0.07%  '
$ws.Range('D30').Value = '133.56'
$ws.Range('E30').Value = '  +0.94%  '
$ws.Range('D31').Value = '2.392'
$ws.Range('E31').Value = '  -0.43%  '
$ws.Range('D32').Value = '6.789'
$ws.Range('E32').Value = '  -0.71%  '
$ws.Range('D33').Value = '1.777.29'
$ws.Range('E33').Value = '  +0.06%  '
$ws.Range('D34').Value = '0.9712'
$ws.Range('E34').Value = '  +0.92%  '
$ws.Range('D35').Value = '0.07482'
$ws.Range('E35').Value = '  -2.20%  '
$ws.Range('D37').Value = '0.02706'
$ws.Range('E37').Value = '  -0.70%  '
$ws.Range('D38').Value = '0.2506'
$ws.Range('E38').Value = '  -1.63%  '
$ws.Range('D39').Value = '0.08793'
$ws.Range('E39').Value = '  -1.25%  '
$ws.Range('D40').Value = '6.086'
$ws.Range('E40').Value = '  -2.70%  '
$ws.Range('B41').Value = 'TheSandbox'
$ws.Range('C41').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D41').Value = '0.7094'
$ws.Range('E41').Value = '  +0.04%  '
$ws.Range('B42').Value = 'TrustWalletToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D42').Value = '1.358'
$ws.Range('E42').Value = '  -0.78%  '
$ws.Range('D43').Value = '12.46'
$ws.Range('E43').Value = '  -1.51%  '
$ws.Range('D44').Value = '15.40'
$ws.Range('E44').Value = '  -0.38%  '
$ws.Range('D45').Value = '0.6521'
$ws.Range('D46').Value = '2.307'
$ws.Range('E46').Value = '  -0.08%  '
$ws.Range('D47').Value = '4.013'
$ws.Range('D48').Value = '131.94'
$ws.Range('E48').Value = '  -0.02%  '
$ws.Range('D49').Value = '0.07961'
$ws.Range('E49').Value = '  +0.33%  '
$ws.Range('D50').Value = '1.200'
$ws.Range('E50').Value = '  -0.42%  '
$ws.Range('D51').Value = '1.215'
$ws.Range('E51').Value = '  +1.81%  '
